$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the distinguishing data between row 3 and row 4
$cols = @("D", "I", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $r3 = $ws.Range("$col`3")
    $r4 = $ws.Range("$col`4")
    $v3 = $r3.Value2()
    $v4 = $r4.Value2()
    $r3.Value = $v4
    $r4.Value = $v3
}
